$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (r2..r10), columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
    2  = @{ E=3; G=27.03890566666666; H=81.11671699999999; I=0.07096188219033728; J=0.07096188219033729; K=3; M=2.044118333333333; N=6.132354999999999; O=0.1776005292722278; P=0.1776005292722278; Q=55.27072278650387; R=497.4365050785349; S=0.01260286783515737; T=0.01260286783515738 }
    3  = @{ E=3; G=27.03890566666666; H=81.11671699999999; I=0.07096188219033728; J=0.07096188219033729; K=3; M=7.059280333333334; N=21.177841; O=0.6133362746356149; P=0.6133362746356149; Q=190.8752150075552; R=1717.876935067997; S=0.04352349646375286; T=0.04352349646375286 }
    4  = @{ E=3; G=27.03890566666666; H=81.11671699999999; I=0.07096188219033728; J=0.07096188219033729; K=3; M=2.406242333333333; N=7.218726999999999; O=0.2090631960921573; P=0.2090631960921573; Q=65.06215946213987; R=585.5594351592589; S=0.01483551789142704; T=0.01483551789142705 }
    5  = @{ E=3; G=345.566579; H=1036.699737; I=0.9069174311350353; J=0.9069174311350354; K=3; M=2.044118333333333; N=6.132354999999999; O=0.1776005292722278; P=0.1776005292722278; Q=706.3789795211816; R=6357.410815690633; S=0.1610690157757915; T=0.1610690157757915 }
    6  = @{ E=3; G=345.566579; H=1036.699737; I=0.9069174311350353; J=0.9069174311350354; K=3; M=7.059280333333334; N=21.177841; O=0.6133362746356149; P=0.6133362746356149; Q=2439.45135499198; R=21955.06219492782; S=0.5562453586144643; T=0.5562453586144644 }
    7  = @{ E=3; G=345.566579; H=1036.699737; I=0.9069174311350353; J=0.9069174311350354; K=3; M=2.406242333333333; N=7.218726999999999; O=0.2090631960921573; P=0.2090631960921573; Q=831.5169313749776; R=7483.652382374798; S=0.1896030567447794; T=0.1896030567447795 }
    8  = @{ E=3; G=8.428738666666668; H=25.286216; I=0.0221206866746274; J=0.02212068667462741; K=3; M=2.044118333333333; N=6.132354999999999; O=0.1776005292722278; P=0.1776005292722278; Q=17.22933923540889; R=155.06405311868; S=0.003928645661278943; T=0.003928645661278944 }
    9  = @{ E=3; G=8.428738666666668; H=25.286216; I=0.0221206866746274; J=0.02212068667462741; K=3; M=7.059280333333334; N=21.177841; O=0.6133362746356149; P=0.6133362746356149; Q=59.50082910440623; R=535.5074619396561; S=0.01356741955739766; T=0.01356741955739766 }
    10 = @{ E=3; G=8.428738666666668; H=25.286216; I=0.0221206866746274; J=0.02212068667462741; K=3; M=2.406242333333333; N=7.218726999999999; O=0.2090631960921573; P=0.2090631960921573; Q=20.28158779633689; R=182.534290167032; S=0.004624621455950799; T=0.0046246214559508 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
